$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.771.50'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '1.864.86'
$ws.Range("E3").Value = '  -1.52%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.38%  '
$ws.Range("D5").Value = '0.7323'
$ws.Range("E5").Value = '  -5.34%  '
$ws.Range("D6").Value = '241.49'
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").Value = '0.3088'
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("D9").Value = '24.56'
$ws.Range("E9").Value = '  -4.54%  '
$ws.Range("D10").Value = '0.07022'
$ws.Range("E10").Value = '  -4.54%  '
$ws.Range("D11").Value = '0.08430'
$ws.Range("E11").Value = '  +4.43%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.902.78'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.7484'
$ws.Range("E13").Value = '  -3.04%  '
$ws.Range("D14").Value = '5.300'
$ws.Range("E14").Value = '  -3.74%  '
$ws.Range("D15").Value = '91.94'
$ws.Range("E15").Value = '  -2.51%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '29.837.36'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = '6.077'
$ws.Range("E17").Value = '  -2.37%  '
$ws.Range("D18").Value = '13.49'
$ws.Range("E18").Value = '  -3.73%  '
$ws.Range("D19").Value = '240.26'
$ws.Range("E19").Value = '  -2.55%  '
$ws.Range("D20").Value = '0.000007745'
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '0.9992'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.137.44'
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("D23").Value = '7.904'
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '0.1560'
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("D26").Value = '9.249'
$ws.Range("E26").Value = '  -2.15%  '
$ws.Range("D27").Value = '162.10'
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").Value = '18.47'
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("D29").Value = '2.007'
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("D30").Value = '1.473'
$ws.Range("E30").Value = '  +3.38%  '
$ws.Range("D31").Value = '1.529'
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("D32").Value = '4.433'
$ws.Range("E32").Value = '  -0.89%  '
$ws.Range("D33").Value = '4.115'
$ws.Range("E33").Value = '  +1.18%  '
$ws.Range("D34").Value = '0.05342'
$ws.Range("E34").Value = '  -4.24%  '
$ws.Range("D35").Value = '1.223'
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("D36").Value = '0.7419'
$ws.Range("E36").Value = '  -1.63%  '
$ws.Range("E37").Value = '  +0.61%  '
$ws.Range("D38").Value = '2.698'
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("D39").Value = '0.01928'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '2.760'
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("D41").Value = '0.4409'
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.091.32'
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '6.040'
$ws.Range("E43").Value = '  +0.47%  '
$ws.Range("D44").Value = '71.70'
$ws.Range("E44").Value = '  -3.89%  '
$ws.Range("D45").Value = '0.8637'
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("D46").Value = '1.002'
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").Value = '101.78'
$ws.Range("E47").Value = '  -0.63%  '
$ws.Range("D48").Value = '7.627'
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("D49").Value = '1.827'
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("D50").Value = '3.001'
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").Value = '2.039.93'
$ws.Range("E51").Value = '  -0.43%  '
